$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text: was "4.3.1.1. Youth education by gender", becomes
#     "4.3.1.1 Youth education by gender" (drop the stray period). This
#     orphans the old shared string (auto-pruned) and appends a new one.
$ws.Range("C1").Value = "4.3.1.1 Youth education by gender"

# --- Add the new 2021 column (M), copying formatting from the matching
#     row's existing K-column cell (same visual style per row), then
#     filling in the 2021 figures.

# Row 2 (blank separator row under the header) just needs the thin-border
# style extended into M.
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

# Row 3: year headers. Also normalizes L3's stray one-off style to match
# K3/M3 (all three end up sharing the same "year header" style).
$ws.Range("K3").Copy()
$ws.Range("L3:M3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2020
$ws.Range("M3").Value = 2021

# Row 4
$ws.Range("K4").Copy()
$ws.Range("L4:M4").PasteSpecial(-4122)
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 10.8

# Row 5
$ws.Range("K5").Copy()
$ws.Range("L5:M5").PasteSpecial(-4122)
$ws.Range("L5").Value = 6.4
$ws.Range("M5").Value = 5.2

# Row 6
$ws.Range("K6").Copy()
$ws.Range("L6:M6").PasteSpecial(-4122)
$ws.Range("L6").Value = 13.5
$ws.Range("M6").Value = 16.2

# Row 7
$ws.Range("K7").Copy()
$ws.Range("L7:M7").PasteSpecial(-4122)
$ws.Range("L7").Value = 24.3
$ws.Range("M7").Value = 24.2

# Row 8
$ws.Range("K8").Copy()
$ws.Range("L8:M8").PasteSpecial(-4122)
$ws.Range("L8").Value = 27.8
$ws.Range("M8").Value = 27.6

# Row 9
$ws.Range("K9").Copy()
$ws.Range("L9:M9").PasteSpecial(-4122)
$ws.Range("L9").Value = 20.9
$ws.Range("M9").Value = 20.9

# Row 10
$ws.Range("K10").Copy()
$ws.Range("L10:M10").PasteSpecial(-4122)
$ws.Range("L10").Value = 26.7
$ws.Range("M10").Value = 28.5

# Row 11
$ws.Range("K11").Copy()
$ws.Range("L11:M11").PasteSpecial(-4122)
$ws.Range("L11").Value = 28.4
$ws.Range("M11").Value = 29.7

# Row 12
$ws.Range("K12").Copy()
$ws.Range("L12:M12").PasteSpecial(-4122)
$ws.Range("L12").Value = 25
$ws.Range("M12").Value = 27.5

# --- Match the saved selection cursor position.
$ws.Range("O2").Select()
